$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(3)
$tf2 = $sh.TextFrame2
$tr2 = $tf2.TextRange
Write-Host "tr2 text: $($tr2.Text)"
$tr2.Text = "9/11/2019"
Write-Host "tr2 text after: $($tr2.Text)"
